# Applies a weekly re-shuffle of the data rows (rows 2-11) in the single
# worksheet of the workbook. Only the "data" columns (D, L, M, N, O, P, Q,
# R, S, T) move between rows; columns A, B, C, E, F, G, H, I, J, K stay put
# because they are identical for every row already.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: target row -> source row (i.e. target row ends up holding the
# values that used to live in the source row).
$rowMap = @{
    2  = 11
    3  = 10
    4  = 6
    5  = 2
    6  = 8
    7  = 9
    8  = 5
    9  = 4
    10 = 3
    11 = 7
}

# Snapshot the "before" values for the columns that move, for every row,
# before we start overwriting anything.
$cols = @("D", "L", "M", "N", "O", "P", "Q", "R", "S", "T")

$snapshot = @{}
foreach ($r in 2..11) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

foreach ($targetRow in $rowMap.Keys) {
    $sourceRow = $rowMap[$targetRow]
    $sourceVals = $snapshot[$sourceRow]
    foreach ($c in $cols) {
        $ws.Range("$c$targetRow").Value = $sourceVals[$c]
    }
}
